$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "loads" sheet: add "v_nom_kv" and "s_base_mva" columns (B, C),
# shifting the former v_nom_pu/p_nom_mw/q_nom_mvar/bus_idx columns
# right into D/E/F/G, and append two new trailing columns (H, I)
# for g_shunt_pu / b_shunt_pu.
# New layout: name | v_nom_kv | s_base_mva | v_nom_pu | p_nom_mw | q_nom_mvar | bus_idx | g_shunt_pu | b_shunt_pu
# ---------------------------------------------------------------
$loads = $wb.Worksheets.Item("loads")

# Header row
$loads.Range("A1").Value = "name"
$loads.Range("B1").Value = "v_nom_kv"
$loads.Range("C1").Value = "s_base_mva"
$loads.Range("D1").Value = "v_nom_pu"
$loads.Range("E1").Value = "p_nom_mw"
$loads.Range("F1").Value = "q_nom_mvar"
$loads.Range("G1").Value = "bus_idx"
$loads.Range("H1").Value = "g_shunt_pu"
$loads.Range("I1").Value = "b_shunt_pu"

# Row 2 - Load 1
$loads.Range("A2").Value = "Load 1"
$loads.Range("B2").Value = 22
$loads.Range("C2").Value = 100
$loads.Range("D2").Value = 1
$loads.Range("E2").Value = 50
$loads.Range("F2").Value = 20
$loads.Range("G2").Value = 1
$loads.Range("H2").Value = 0
$loads.Range("I2").Value = 0

# Row 3 - Load 2
$loads.Range("A3").Value = "Load 2"
$loads.Range("B3").Value = 132
$loads.Range("C3").Value = 100
$loads.Range("D3").Value = 1
$loads.Range("E3").Value = 5
$loads.Range("F3").Value = 2
$loads.Range("G3").Value = 3
$loads.Range("H3").Value = 0
$loads.Range("I3").Value = 0

# Row 4 - Load 3
$loads.Range("A4").Value = "Load 3"
$loads.Range("B4").Value = 132
$loads.Range("C4").Value = 100
$loads.Range("D4").Value = 1
$loads.Range("E4").Value = 5
$loads.Range("F4").Value = 2
$loads.Range("G4").Value = 4
$loads.Range("H4").Value = 0
$loads.Range("I4").Value = 0

[void]$loads.Range("I5").Select()

# ---------------------------------------------------------------
# "trafos" sheet: insert "v_base_kV" column (E) after "V_lv_kV",
# shifting the old V_SCH_pu..tap_max columns right by one (F..O).
# New layout: name | S_nom | V_hv_kV | V_lv_kV | v_base_kV | V_SCH_pu | P_Cu_pu | I_E_pu | P_Fe_pu | idx_hv | idx_lv | tap_pos | tap_change | tap_min | tap_max
# ---------------------------------------------------------------
$trafos = $wb.Worksheets.Item("trafos")

# Header row
$trafos.Range("A1").Value = "name"
$trafos.Range("B1").Value = "S_nom"
$trafos.Range("C1").Value = "V_hv_kV"
$trafos.Range("D1").Value = "V_lv_kV"
$trafos.Range("E1").Value = "v_base_kV"
$trafos.Range("F1").Value = "V_SCH_pu"
$trafos.Range("G1").Value = "P_Cu_pu"
$trafos.Range("H1").Value = "I_E_pu"
$trafos.Range("I1").Value = "P_Fe_pu"
$trafos.Range("J1").Value = "idx_hv"
$trafos.Range("K1").Value = "idx_lv"
$trafos.Range("L1").Value = "tap_pos"
$trafos.Range("M1").Value = "tap_change"
$trafos.Range("N1").Value = "tap_min"
$trafos.Range("O1").Value = "tap_max"

# Row 2 - T1
$trafos.Range("A2").Value = "T1"
$trafos.Range("B2").Value = 20
$trafos.Range("C2").Value = 132
$trafos.Range("D2").Value = 22
$trafos.Range("E2").Value = 132
$trafos.Range("F2").Value = 0.03
$trafos.Range("G2").Value = 0.005
$trafos.Range("H2").Value = 0.03
$trafos.Range("I2").Value = 0.001
$trafos.Range("J2").Value = 1
$trafos.Range("K2").Value = 3
$trafos.Range("L2").Value = 0
$trafos.Range("M2").Value = 0.01
$trafos.Range("N2").Value = -5
$trafos.Range("O2").Value = 5

[void]$trafos.Range("E2").Select()

$wb.Save()
